# Correção das notas do fórum para matc65 em 2021.2
# Zera as colunas de visualizações diárias (B:H), total_views (I) e nota_view (J)
# para todas as linhas de alunos (linhas 2 a 50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:J50").Value = 0
